# Applies the "update pccv, project outline" edit to the project outline
# document: strips the per-task "- <Name>" owner suffixes, updates the
# Week 9/10 schedule paragraphs, and rewords the Week 10/11 heading.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $null = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2)
}

# --- remove / tidy the "- <owner>" suffixes -------------------------------

Replace-Text "Design database schema (28/3 - 3/4) – Duy" "Design database schema (28/3 - 3/4)"
Replace-Text "SRS: Hoàng + Tấn + Dũng" "SRS"
Replace-Text "CRUD – Duy" "CRUD "
Replace-Text "Filter - Duy" "Filter "
Replace-Text "Sort - Hoàng" "Sort "
Replace-Text "Limit fileds - Dũng" "Limit fileds "
Replace-Text "Pagination - Tấn" "Pagination "
Replace-Text "Login - Hoàng" "Login "
Replace-Text "Sign up - Tấn" "Sign up "
Replace-Text "Forgot password - Duy" "Forgot password "
Replace-Text "Reset password - Duy" "Reset password "
Replace-Text "Update password - Tấn" "Update password "
Replace-Text "Update user info - Hoàng" "Update user info "
Replace-Text "Delete user - Dũng" "Delete user "
Replace-Text "Protect route - Duy" "Protect route "
Replace-Text "Authorization - Duy" "Authorization "
Replace-Text "API docs - Duy" "API docs "

# This lone "- Duy" belongs to the "CRUD" bullet under "Review"; every other
# "- <owner>" phrase above has already been neutralised, so this is now the
# only remaining match in the whole document.
Replace-Text " - Duy" " "

# --- Week 9 / Week 10 schedule rewording ----------------------------------

Replace-Text "Week 9: Render data into template (16/5 - 22/5) (Duy)" "Week 9, 10: Render data into template (16/5 - 29/5) "

Replace-Text "Week 10: Payments, Email, File uploads (23/5 - 29/5)" "Week 11: Payments, Email, File uploads (30/5 – 5/6)"
